$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 1218
$ws.Range("I31").Value = 796.6667
$ws.Range("J31").Value = 1850
$ws.Range("K31").Value = 2390.0001
$ws.Range("L31").Value = 5550
$ws.Range("M31").Value = -2160.0001
$ws.Range("N31").Value = -6010
$ws.Range("H40").Value = 1104.3636
$ws.Range("I40").Value = 885.5714
$ws.Range("J40").Value = 1487.25
$ws.Range("K40").Value = 885.5714
$ws.Range("L40").Value = 1487.25
$ws.Range("M40").Value = -710.5714
$ws.Range("N40").Value = -1837.25
$ws.Range("H69").Value = 4676
$ws.Range("I69").Value = 4542
$ws.Range("J69").Value = 4810
$ws.Range("K69").Value = 13626
$ws.Range("L69").Value = 14430
$ws.Range("M69").Value = -12752
$ws.Range("N69").Value = -16178
$ws.Range("H72").Value = 4676
$ws.Range("I72").Value = 4542
$ws.Range("J72").Value = 4810
$ws.Range("K72").Value = 40878
$ws.Range("L72").Value = 43290
$ws.Range("M72").Value = -36510
$ws.Range("N72").Value = -52026
$ws.Range("H98").Value = 2229.9
$ws.Range("I98").Value = 1247.2142
$ws.Range("J98").Value = 4522.8335
$ws.Range("K98").Value = 1247.2142
$ws.Range("L98").Value = 4522.8335
$ws.Range("M98").Value = 250.7858000000001
$ws.Range("N98").Value = -7518.8335
$ws.Range("H122").Value = 2229.9
$ws.Range("I122").Value = 1247.2142
$ws.Range("J122").Value = 4522.8335
$ws.Range("K122").Value = 3741.6426
$ws.Range("L122").Value = 13568.5005
$ws.Range("M122").Value = -1291.6426
$ws.Range("N122").Value = -18468.5005
$ws.Range("H136").Value = 31082.857
$ws.Range("J136").Value = 31082.857
$ws.Range("L136").Value = 31082.857
$ws.Range("N136").Value = -41282.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5817.404
$ws.Range("I32").Value = 4604.148
$ws.Range("J32").Value = 11277.056
$ws.Range("K32").Value = 4604.148
$ws.Range("L32").Value = 11277.056
$ws.Range("M32").Value = -4317.148
$ws.Range("N32").Value = -11851.056
$ws.Range("H74").Value = 566
$ws.Range("I74").Value = 556.9
$ws.Range("K74").Value = 556.9
$ws.Range("M74").Value = 317.1
$ws.Range("H77").Value = 566
$ws.Range("I77").Value = 556.9
$ws.Range("K77").Value = 2784.5
$ws.Range("M77").Value = 1583.5
$ws.Range("H132").Value = 2563.1538
$ws.Range("I132").Value = 2041.0358
$ws.Range("J132").Value = 3892.182
$ws.Range("K132").Value = 6123.107400000001
$ws.Range("L132").Value = 11676.546
$ws.Range("M132").Value = -3593.107400000001
$ws.Range("N132").Value = -16736.546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 557
$ws.Range("I5").Value = 463.33334
$ws.Range("J5").Value = 697.5
$ws.Range("K5").Value = 463.33334
$ws.Range("L5").Value = 697.5
$ws.Range("M5").Value = -350.33334
$ws.Range("N5").Value = -923.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 44501.332
$ws.Range("I3").Value = 19334
$ws.Range("J3").Value = 69668.664
$ws.Range("K3").Value = 19334
$ws.Range("L3").Value = 69668.664
$ws.Range("M3").Value = -19221
$ws.Range("N3").Value = -69894.664
$ws.Range("H31").Value = 2225234.5
$ws.Range("I31").Value = 3573266
$ws.Range("J31").Value = 4947.1763
$ws.Range("K31").Value = 3573266
$ws.Range("L31").Value = 4947.1763
$ws.Range("M31").Value = -3572971
$ws.Range("N31").Value = -5537.1763
$ws.Range("H34").Value = 2225234.5
$ws.Range("I34").Value = 3573266
$ws.Range("J34").Value = 4947.1763
$ws.Range("K34").Value = 3573266
$ws.Range("L34").Value = 4947.1763
$ws.Range("M34").Value = -3573064
$ws.Range("N34").Value = -5351.1763
$ws.Range("H86").Value = 8374.875
$ws.Range("I86").Value = 8500
$ws.Range("J86").Value = 8333.166999999999
$ws.Range("K86").Value = 8500
$ws.Range("L86").Value = 8333.166999999999
$ws.Range("M86").Value = -7377
$ws.Range("N86").Value = -10579.167
$ws.Range("H89").Value = 8374.875
$ws.Range("I89").Value = 8500
$ws.Range("J89").Value = 8333.166999999999
$ws.Range("K89").Value = 42500
$ws.Range("L89").Value = 41665.835
$ws.Range("M89").Value = -36884
$ws.Range("N89").Value = -52897.835
$ws.Range("H122").Value = 3678.4583
$ws.Range("I122").Value = 2814.5715
$ws.Range("J122").Value = 4034.1765
$ws.Range("K122").Value = 8443.7145
$ws.Range("L122").Value = 12102.5295
$ws.Range("M122").Value = -5993.7145
$ws.Range("N122").Value = -17002.5295
$ws.Range("H132").Value = 1600.1
$ws.Range("I132").Value = 1128.2325
$ws.Range("J132").Value = 4498.7144
$ws.Range("K132").Value = 3384.6975
$ws.Range("L132").Value = 13496.1432
$ws.Range("M132").Value = -854.6975000000002
$ws.Range("N132").Value = -18556.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 16257.571
$ws.Range("J41").Value = 16257.571
$ws.Range("L41").Value = 48772.713
$ws.Range("N41").Value = -49448.713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 62715.715
$ws.Range("J4").Value = 62715.715
$ws.Range("L4").Value = 62715.715
$ws.Range("N4").Value = -62939.715
$ws.Range("H80").Value = 2725
$ws.Range("I80").Value = 2740
$ws.Range("J80").Value = 2700
$ws.Range("K80").Value = 2740
$ws.Range("L80").Value = 2700
$ws.Range("M80").Value = -1742
$ws.Range("N80").Value = -4696
$ws.Range("H83").Value = 2725
$ws.Range("I83").Value = 2740
$ws.Range("J83").Value = 2700
$ws.Range("K83").Value = 13700
$ws.Range("L83").Value = 13500
$ws.Range("M83").Value = -8708
$ws.Range("N83").Value = -23484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 125001816
$ws.Range("I22").Value = 250000400
$ws.Range("J22").Value = 3225
$ws.Range("K22").Value = 250000400
$ws.Range("L22").Value = 3225
$ws.Range("M22").Value = -250000105
$ws.Range("N22").Value = -3815
$ws.Range("H27").Value = 125001816
$ws.Range("I27").Value = 250000400
$ws.Range("J27").Value = 3225
$ws.Range("K27").Value = 250000400
$ws.Range("L27").Value = 3225
$ws.Range("M27").Value = -250000293
$ws.Range("N27").Value = -3439
$ws.Range("H55").Value = 1235.6364
$ws.Range("I55").Value = 157.14285
$ws.Range("J55").Value = 3123
$ws.Range("K55").Value = 157.14285
$ws.Range("L55").Value = 3123
$ws.Range("M55").Value = 15.85714999999999
$ws.Range("N55").Value = -3469
$ws.Range("H82").Value = 2254.7058
$ws.Range("I82").Value = 1760
$ws.Range("J82").Value = 2961.4285
$ws.Range("K82").Value = 1760
$ws.Range("L82").Value = 2961.4285
$ws.Range("M82").Value = -1399
$ws.Range("N82").Value = -3683.4285
$ws.Range("H85").Value = 2254.7058
$ws.Range("I85").Value = 1760
$ws.Range("J85").Value = 2961.4285
$ws.Range("K85").Value = 1760
$ws.Range("L85").Value = 2961.4285
$ws.Range("M85").Value = -512
$ws.Range("N85").Value = -5457.4285
$ws.Range("H132").Value = 2175.6943
$ws.Range("I132").Value = 1446.875
$ws.Range("K132").Value = 4340.625
$ws.Range("M132").Value = -1810.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 857.7778
$ws.Range("I81").Value = 840
$ws.Range("K81").Value = 1680
$ws.Range("M81").Value = -619
$ws.Range("H84").Value = 857.7778
$ws.Range("I84").Value = 840
$ws.Range("K84").Value = 8400
$ws.Range("M84").Value = -3096
$ws.Range("H122").Value = 669346.1
$ws.Range("I122").Value = 835249.4399999999
$ws.Range("K122").Value = 2505748.32
$ws.Range("M122").Value = -2503298.32
$ws.Range("H132").Value = 193081.86
$ws.Range("I132").Value = 234646.33
$ws.Range("K132").Value = 703938.99
$ws.Range("M132").Value = -701408.99
$ws.Range("H136").Value = 732.1579
$ws.Range("I136").Value = 341.14285
$ws.Range("J136").Value = 1827
$ws.Range("K136").Value = 1023.42855
$ws.Range("L136").Value = 5481
$ws.Range("M136").Value = 1526.57145
$ws.Range("N136").Value = -10581

Write-Output "done"